# update function to add row color based on value.
# Replace the conditional-formatting-based row highlighting with literal
# values + a directly-applied fill color per row (computed "based on value"),
# and duplicate the highlighted summary row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$highlight = 15658734   # RGB(238,238,238) == FFEEEEEE

# --- Row 1: plain numeric header row (no longer shared-string text) ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3
# Touch the alignment so this row gets its own (non-default) style slot,
# matching the "normal" look used again on row 3.
$ws.Range("A1:C1").IndentLevel = 0

# --- Row 2: "a/b/c" - value row that gets the highlight color ---
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "b"
$ws.Range("C2").Value = "c"
$ws.Range("A2:C2").Interior.Color = $highlight

# --- Row 3: "d/e/f/x" - normal (unhighlighted) row ---
$ws.Range("A3").Value = "d"
$ws.Range("B3").Value = "e"
$ws.Range("C3").Value = "f"
$ws.Range("D3").Value = "x"
$ws.Range("A3:D3").IndentLevel = 0

# --- Row 4: "c/x1" - new highlighted summary row ---
$ws.Range("A4").Value = "c"
$ws.Range("B4").Value = "x1"
$ws.Range("A4:B4").Interior.Color = $highlight

# --- Row 5: "c/x1" - duplicate highlighted summary row ---
$ws.Range("A5").Value = "c"
$ws.Range("B5").Value = "x1"
$ws.Range("A5:B5").Interior.Color = $highlight

# The old dxf-based conditional formatting is superseded by the literal
# fills applied above, so drop it.
$ws.Range("A4:E4").FormatConditions.Delete()

Write-Output "done"
